$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("catalogo")

# Grow the "Tabla1" structured table by one row (A1:I3 -> A1:I4),
# which also expands the table's AutoFilter range.
$lo = $ws.ListObjects.Item("Tabla1")
$newRow = $lo.ListRows.Add()

# Fill in the new product row (4): nombre, tipo, categoría, precio, imagen1
$ws.Range("A4").Value = "Kioto Dorado"
$ws.Range("B4").Value = "bolsos"
$ws.Range("C4").Value = "Kioto"
$ws.Range("D4").Value = 25
$ws.Range("D4").NumberFormat = $ws.Range("D3").NumberFormat
$ws.Range("E4").Value = "imagen3_1.jpg"

# Extend the "tipo" and "categoría" list-validation ranges to cover the new row
$ws.Range("B2:B4").Validation.Delete()
$ws.Range("B2:B4").Validation.Add(3, 1, 1, "tipos")

$ws.Range("C2:C4").Validation.Delete()
$ws.Range("C2:C4").Validation.Add(3, 1, 1, "categoria")

# Move the active selection to reflect where editing left off
$null = $ws.Range("F4").Select()
